$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old scratch rows (20 and 26)
$ws.Range("E20:F20").ClearContents()
$ws.Range("E26").ClearContents()

# Add the new value
$ws.Range("C22").Value = "拉之前的"

# Update selection to match the final state
$ws.Range("C22").Select()
